$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuille5")

# Insert two new columns before the existing "Temps utilisé" column (E),
# pushing it to G, to make room for "Prénom" (E) and "Note/10,00" (F).
$ws.Columns("E:F").Insert()

# Headers
$ws.Range("E1").Value = "Prénom"
$ws.Range("F1").Value = "Note/10,00"

# Data: column E duplicates the first name already present in column B,
# column F holds the grade out of 10.
$firstNames = @{
    2  = "Houzefa"
    3  = "Yasmine"
    4  = "Zina"
    5  = "Aboubaker"
    6  = "Yasmine"
    7  = "Hassan Mahamat"
    8  = "Yacine"
    9  = "Paola"
    10 = "Rodolphe"
    11 = "Nouh"
    12 = "Iness"
    13 = "Zakaria"
    14 = "Christian"
}

$grades = @{
    2  = "7,83"
    3  = "7,83"
    4  = "7,28"
    5  = "7,98"
    6  = "7,83"
    7  = "7,52"
    8  = "8,07"
    9  = "7,16"
    10 = "8,88"
    11 = "7,35"
    12 = "6,51"
    13 = "7,70"
    14 = "6,72"
}

foreach ($r in 2..14) {
    $ws.Cells.Item($r, 5).Value = $firstNames[$r]
    $ws.Cells.Item($r, 6).Value = $grades[$r]
}

# Match the highlighted-row styling already used in column B for rows 10 and 14.
$ws.Range("B10").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("B14").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# Restore the values overwritten by the format paste.
$ws.Cells.Item(10, 5).Value = $firstNames[10]
$ws.Cells.Item(14, 5).Value = $firstNames[14]

Write-Host "done"
